# Update Leve profit data cells across multiple sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# Each block targets one row identified by its Leve entry; values come from a scheduled data refresh.
$wb = $excel.ActiveWorkbook

# ALC row 98
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 3010.7368
$ws.Range("I98").Value = 2970.818
$ws.Range("K98").Value = 2970.818
$ws.Range("M98").Value = -1472.818

# ALC row 99
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H99").Value = 1386.6666
$ws.Range("I99").Value = 930.75
$ws.Range("J99").Value = 2298.5
$ws.Range("K99").Value = 2792.25
$ws.Range("L99").Value = 6895.5
$ws.Range("M99").Value = -1294.25
$ws.Range("N99").Value = -9891.5

# ALC row 104
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H104").Value = 225.33333
$ws.Range("I104").Value = 225.33333
$ws.Range("K104").Value = 675.99999
$ws.Range("M104").Value = 1071.00001

# ALC row 122
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 3010.7368
$ws.Range("I122").Value = 2970.818
$ws.Range("K122").Value = 8912.454000000002
$ws.Range("M122").Value = -6462.454000000002

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1462.2
$ws.Range("I137").Value = 1322.6875
$ws.Range("J137").Value = 2020.25
$ws.Range("K137").Value = 3968.0625
$ws.Range("L137").Value = 6060.75
$ws.Range("M137").Value = -1418.0625
$ws.Range("N137").Value = -11160.75

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 1580.7778
$ws.Range("I138").Value = 454.42856
$ws.Range("J138").Value = 1975
$ws.Range("K138").Value = 1363.28568
$ws.Range("L138").Value = 5925
$ws.Range("M138").Value = 3776.71432
$ws.Range("N138").Value = -16205

# ARM row 13
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H13").Value = 17914.666
$ws.Range("J13").Value = 17914.666
$ws.Range("L13").Value = 17914.666
$ws.Range("N13").Value = -18202.666

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1533.2727
$ws.Range("I74").Value = 1533.2727
$ws.Range("K74").Value = 1533.2727
$ws.Range("M74").Value = -659.2727

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1533.2727
$ws.Range("I77").Value = 1533.2727
$ws.Range("K77").Value = 7666.363499999999
$ws.Range("M77").Value = -3298.363499999999

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1367.4828
$ws.Range("I132").Value = 1318.8
$ws.Range("K132").Value = 3956.4
$ws.Range("M132").Value = -1426.4

# BSM row 60
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H60").Value = 41498.5
$ws.Range("J60").Value = 41498.5
$ws.Range("L60").Value = 41498.5
$ws.Range("N60").Value = -42696.5

# BSM row 94
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1708
$ws.Range("I94").Value = 1995.6666
$ws.Range("J94").Value = 845
$ws.Range("K94").Value = 1995.6666
$ws.Range("L94").Value = 845
$ws.Range("M94").Value = -1544.6666
$ws.Range("N94").Value = -1747

# BSM row 99
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2333.1667
$ws.Range("I99").Value = 2070
$ws.Range("K99").Value = 2070
$ws.Range("M99").Value = -572

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3656.5
$ws.Range("I31").Value = 4612
$ws.Range("J31").Value = 3338
$ws.Range("K31").Value = 4612
$ws.Range("L31").Value = 3338
$ws.Range("M31").Value = -4317
$ws.Range("N31").Value = -3928

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3656.5
$ws.Range("I34").Value = 4612
$ws.Range("J34").Value = 3338
$ws.Range("K34").Value = 4612
$ws.Range("L34").Value = 3338
$ws.Range("M34").Value = -4410
$ws.Range("N34").Value = -3742

# CRP row 36
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H36").Value = 30000
$ws.Range("J36").Value = 30000
$ws.Range("L36").Value = 30000
$ws.Range("N36").Value = -30776

# CRP row 40
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H40").Value = 30000
$ws.Range("J40").Value = 30000
$ws.Range("L40").Value = 30000
$ws.Range("N40").Value = -30320

# CRP row 41
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 4250
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("N41").ClearContents()

# CRP row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2404.7368
$ws.Range("I58").Value = 2317.7144
$ws.Range("J58").Value = 2648.4
$ws.Range("K58").Value = 2317.7144
$ws.Range("L58").Value = 2648.4
$ws.Range("M58").Value = -2114.7144
$ws.Range("N58").Value = -3054.4

# CRP row 59
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 32130.416
$ws.Range("J59").Value = 33678.363
$ws.Range("L59").Value = 33678.363
$ws.Range("N59").Value = -35968.363

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1542.5
$ws.Range("J132").Value = 1312.75
$ws.Range("L132").Value = 3938.25
$ws.Range("N132").Value = -8998.25

# CRP row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2404.7368
$ws.Range("I136").Value = 2317.7144
$ws.Range("J136").Value = 2648.4
$ws.Range("K136").Value = 6953.1432
$ws.Range("L136").Value = 7945.200000000001
$ws.Range("M136").Value = -4403.1432
$ws.Range("N136").Value = -13045.2

# CUL row 2
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 287
$ws.Range("I2").Value = 70
$ws.Range("K2").Value = 420
$ws.Range("M2").Value = -307

# CUL row 109
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 2053.25
$ws.Range("I109").Value = 2308.6667
$ws.Range("K109").Value = 6926.000100000001
$ws.Range("M109").Value = -5886.000100000001

# CUL row 122
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1326.3334
$ws.Range("I122").Value = 990
$ws.Range("J122").Value = 1494.5
$ws.Range("K122").Value = 8910
$ws.Range("L122").Value = 13450.5
$ws.Range("M122").Value = -6460
$ws.Range("N122").Value = -18350.5

# CUL row 130
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H130").Value = 0
$ws.Range("I130").Value = 0
$ws.Range("K130").Value = 0
$ws.Range("M130").ClearContents()

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 2974.5
$ws.Range("I131").Value = 1598
$ws.Range("K131").Value = 4794
$ws.Range("M131").Value = 246

# GSM row 10
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 3
$ws.Range("I10").Value = 3
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 166

# GSM row 43
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 36444.25
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 36444.25
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 36444.25
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -36746.25

# LTW row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3346.077
$ws.Range("I46").Value = 2250.25
$ws.Range("K46").Value = 2250.25
$ws.Range("M46").Value = -2062.25

# LTW row 61
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 5519.8
$ws.Range("I61").Value = 5519.8
$ws.Range("K61").Value = 5519.8
$ws.Range("M61").Value = -5317.8

# LTW row 113
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 5519.8
$ws.Range("I113").Value = 5519.8
$ws.Range("K113").Value = 5519.8
$ws.Range("M113").Value = -3349.8

# WVR row 41
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 15883.7
$ws.Range("I41").Value = 16542
$ws.Range("K41").Value = 16542
$ws.Range("M41").Value = -16152
